# Add new testcase rows to "Test cases" sheet and tweak a new sql-script
# related entry, per commit "Add testcase and new sql script".

$wb = $excel.ActiveWorkbook
$wsTest = $wb.Worksheets.Item("Test cases")
$wsTodo = $wb.Worksheets.Item("TODOs")

# --- New testcase rows (sparse, every other row, mirroring the author's layout) ---

# Shared-string creation order matches the author's actual typing order
# (it is not always a simple row-by-row A-then-D sweep: the "Admin ..."
# block below was filled by pasting the Expected-Result column first and
# the Description column second).

$wsTest.Range("A30").Value = "Test logout when cart have item"
$wsTest.Range("D30").Value = "Cart still have item"

$wsTest.Range("A32").Value = "Readd an item into cart"
$wsTest.Range("D32").Value = "Warning pop up"

$wsTest.Range("D34").Value = "See newly added user"
$wsTest.Range("D36").Value = "See error and can not add"
$wsTest.Range("D38").Value = "See newly updated information"
$wsTest.Range("D40").Value = "See error and can not update user"

$wsTest.Range("A34").Value = "Admin add user"
$wsTest.Range("A36").Value = "Admin add user empty field"
$wsTest.Range("A38").Value = "Admin update user information"
$wsTest.Range("A40").Value = "Admin update user empty field"

$wsTest.Range("A42").Value = "Staff/Admin add product"
$wsTest.Range("D42").Value = "See newly added product"

$wsTest.Range("A44").Value = "Staff/Admin add product missing fields"
$wsTest.Range("D44").Value = "See error and can not add new product"

$wsTest.Range("A46").Value = "Staff/Admin update product information"
$wsTest.Range("D46").Value = "See newly updated information"
$wsTest.Range("A46:D47").HorizontalAlignment = -4108
$wsTest.Range("A46:D47").VerticalAlignment = -4108

$wsTest.Range("A48").Value = "Staff/Admin add product missing fields"
$wsTest.Range("D48").Value = "See error and can not update product"
$wsTest.Range("A48:D48").HorizontalAlignment = -4108
$wsTest.Range("A48:D48").VerticalAlignment = -4108

$wsTest.Range("A50").Value = "Staff/Admin update product image"
$wsTest.Range("D50").Value = "See newly updated product image"

$wsTest.Range("A52").Value = "Staff/Admin update wrong/empty image file"
$wsTest.Range("D52").Value = "Cannot update"

$wsTest.Range("A54").Value = "View category, sort, filter"

# --- View / selection state ---

$wsTest.Activate()
$excel.ActiveWindow.Zoom = 115
$wsTest.Range("A56").Select()

$wsTodo.Activate()
$wsTodo.Range("B17").Select()

$wsTest.Activate()
